$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-26 22:18:21'
$ws.Range('O2').Value = '5.5 °C'
$ws.Range('E3').Value = '2026-02-26 22:18:24'
$ws.Range('L3').Value = '21.2 km/h - 121º 21:36 TU'
$ws.Range('E4').Value = '2026-02-26 22:18:26'
$ws.Range('O4').Value = '10.5 °C'
$ws.Range('E5').Value = '2026-02-26 22:18:29'
$ws.Range('E6').Value = '2026-02-26 22:18:31'
$ws.Range('O6').Value = '11.7 °C'
$ws.Range('E7').Value = '2026-02-26 22:18:34'
$ws.Range('J7').Value = '1027.3 hPa'
$ws.Range('K7').Value = '15.4 MJ/m2'
$ws.Range('E8').Value = '2026-02-26 22:18:36'
$ws.Range('O8').Value = '11.1 °C'
$ws.Range('E9').Value = '2026-02-26 22:18:39'
$ws.Range('H9').Value = '''87%'
$ws.Range('O9').Value = '11.9 °C'
$ws.Range('E10').Value = '2026-02-26 22:18:41'
$ws.Range('O10').Value = '9.3 °C'
$ws.Range('E11').Value = '2026-02-26 22:18:44'
$ws.Range('H11').Value = '''69%'
$ws.Range('O11').Value = '8.5 °C'
$ws.Range('E12').Value = '2026-02-26 22:18:46'
$ws.Range('N12').Value = '7.1 °C 21:59 TU'
$ws.Range('O12').Value = '11.2 °C'
$ws.Range('E13').Value = '2026-02-26 22:18:48'
$ws.Range('G13').Value = '3 cm'
$ws.Range('J13').Value = '1028.3 hPa'
$ws.Range('O13').Value = '7.1 °C'
$ws.Range('E14').Value = '2026-02-26 22:18:51'
$ws.Range('O14').Value = '11.2 °C'
$ws.Range('E15').Value = '2026-02-26 22:18:53'
$ws.Range('H15').Value = '''85%'
$ws.Range('O15').Value = '11.4 °C'
$ws.Range('E16').Value = '2026-02-26 22:18:56'
$ws.Range('H16').Value = '''41%'
$ws.Range('E17').Value = '2026-02-26 22:18:58'
$ws.Range('E18').Value = '2026-02-26 22:19:01'
$ws.Range('E19').Value = '2026-02-26 22:19:03'
$ws.Range('O19').Value = '11.3 °C'
$ws.Range('E20').Value = '2026-02-26 22:19:06'
$ws.Range('O20').Value = '2.5 °C'
$ws.Range('E21').Value = '2026-02-26 22:19:08'
$ws.Range('O21').Value = '10.0 °C'
$ws.Range('E22').Value = '2026-02-26 22:19:11'
$ws.Range('H22').Value = '''52%'
$ws.Range('E23').Value = '2026-02-26 22:19:13'
$ws.Range('E24').Value = '2026-02-26 22:19:16'
$ws.Range('H24').Value = '''75%'
$ws.Range('O24').Value = '10.3 °C'
$ws.Range('E25').Value = '2026-02-26 22:19:19'
$ws.Range('E26').Value = '2026-02-26 22:19:22'
$ws.Range('J26').Value = '1024.5 hPa'
$ws.Range('O26').Value = '10.7 °C'
$ws.Range('E27').Value = '2026-02-26 22:19:24'
$ws.Range('L27').Value = '22.0 km/h - 241º 21:39 TU'
$ws.Range('E28').Value = '2026-02-26 22:19:27'
$ws.Range('N28').Value = '5.3 °C 21:59 TU'
$ws.Range('O28').Value = '10.6 °C'
$ws.Range('E29').Value = '2026-02-26 22:19:29'
$ws.Range('H29').Value = '''87%'
$ws.Range('N29').Value = '7.1 °C 21:33 TU'
$ws.Range('O29').Value = '11.4 °C'
$ws.Range('E30').Value = '2026-02-26 22:19:32'
$ws.Range('E31').Value = '2026-02-26 22:19:36'
$ws.Range('E32').Value = '2026-02-26 22:19:38'
$ws.Range('H32').Value = '''68%'
$ws.Range('O32').Value = '7.7 °C'
$ws.Range('E33').Value = '2026-02-26 22:19:41'
$ws.Range('H33').Value = '''54%'
$ws.Range('E34').Value = '2026-02-26 22:19:43'
$ws.Range('H34').Value = '''48%'
$ws.Range('O34').Value = '4.6 °C'
$ws.Range('E35').Value = '2026-02-26 22:19:46'
$ws.Range('K35').Value = '16.7 MJ/m2'
$ws.Range('O35').Value = '12.0 °C'
$ws.Range('E36').Value = '2026-02-26 22:19:48'
$ws.Range('H36').Value = '''87%'
$ws.Range('O36').Value = '12.4 °C'
$ws.Range('E37').Value = '2026-02-26 22:19:50'
$ws.Range('H37').Value = '''75%'
$ws.Range('O37').Value = '7.6 °C'
$ws.Range('E38').Value = '2026-02-26 22:19:53'
$ws.Range('K38').Value = '14.8 MJ/m2'
$ws.Range('O38').Value = '11.0 °C'
$ws.Range('E39').Value = '2026-02-26 22:19:55'
$ws.Range('E40').Value = '2026-02-26 22:19:58'
$ws.Range('J40').Value = '1027.5 hPa'
$ws.Range('O40').Value = '9.4 °C'
$ws.Range('E41').Value = '2026-02-26 22:20:00'
$ws.Range('H41').Value = '''85%'
$ws.Range('J41').Value = '1027.1 hPa'
$ws.Range('E42').Value = '2026-02-26 22:20:03'
$ws.Range('E43').Value = '2026-02-26 22:20:05'
$ws.Range('O43').Value = '9.4 °C'
$ws.Range('E44').Value = '2026-02-26 22:20:08'
$ws.Range('H44').Value = '''55%'
$ws.Range('E45').Value = '2026-02-26 22:20:10'
$ws.Range('E46').Value = '2026-02-26 22:20:12'
$ws.Range('J46').Value = '1027.1 hPa'
